# Apply updated loading_percent results (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 4.917212931948922
$ws.Range("D2").Value = 8.167369351208349
$ws.Range("E2").Value = 11.99738351296328
$ws.Range("F2").Value = 35.16875242850849
$ws.Range("G2").Value = 3.64068724090915
$ws.Range("J2").Value = 9.278083156453967
$ws.Range("M2").Value = 24.88075059864842
$ws.Range("N2").Value = 17.35447637075239
$ws.Range("O2").Value = 27.26808620940093

# Row 3
$ws.Range("C3").Value = 4.750119464859718
$ws.Range("D3").Value = 8.179503095209968
$ws.Range("E3").Value = 12.05463908399377
$ws.Range("F3").Value = 34.9631354881279
$ws.Range("G3").Value = 3.64444024518847
$ws.Range("J3").Value = 9.321114434380858
$ws.Range("M3").Value = 24.17427427824269
$ws.Range("N3").Value = 17.09569665456516
$ws.Range("O3").Value = 27.08644799104477

# Row 4
$ws.Range("C4").Value = 4.646162616357148
$ws.Range("D4").Value = 8.187451942888639
$ws.Range("E4").Value = 12.09152148324506
$ws.Range("F4").Value = 34.84836795875881
$ws.Range("G4").Value = 3.646864997017706
$ws.Range("J4").Value = 9.34877062059031
$ws.Range("M4").Value = 23.73097930278247
$ws.Range("N4").Value = 16.93687677234264
$ws.Range("O4").Value = 26.98332921512337

# Row 5
$ws.Range("C5").Value = 4.603532865019321
$ws.Range("D5").Value = 8.190816853145179
$ws.Range("E5").Value = 12.10698707767452
$ws.Range("F5").Value = 34.80452020706537
$ws.Range("G5").Value = 3.647883491431646
$ws.Range("J5").Value = 9.360352286857824
$ws.Range("M5").Value = 23.54818489355255
$ws.Range("N5").Value = 16.87225031281281
$ws.Range("O5").Value = 26.94345513850023

# Row 6
$ws.Range("C6").Value = 4.596440452698929
$ws.Range("D6").Value = 8.19138319392686
$ws.Range("E6").Value = 12.10958148480595
$ws.Range("F6").Value = 34.79741661514781
$ws.Range("G6").Value = 3.64805445040567
$ws.Range("J6").Value = 9.362294260921404
$ws.Range("M6").Value = 23.51770976550928
$ws.Range("N6").Value = 16.86152704277186
$ws.Range("O6").Value = 26.93696468571127

# Row 7
$ws.Range("C7").Value = 4.64558867147369
$ws.Range("D7").Value = 8.187496813934739
$ws.Range("E7").Value = 12.09172829164696
$ws.Range("F7").Value = 34.84776474566191
$ws.Range("G7").Value = 3.646878609597311
$ws.Range("J7").Value = 9.348925552198029
$ws.Range("M7").Value = 23.72852243531446
$ws.Range("N7").Value = 16.93600471623157
$ws.Range("O7").Value = 26.9827827227267

# Row 8
$ws.Range("C8").Value = 4.859925192702277
$ws.Range("D8").Value = 8.1714497418567
$ws.Range("E8").Value = 12.01676758238604
$ws.Range("F8").Value = 35.09549446808343
$ws.Range("G8").Value = 3.641956358657358
$ws.Range("J8").Value = 9.292664688353154
$ws.Range("M8").Value = 24.63927400375151
$ws.Range("N8").Value = 17.26527943722603
$ws.Range("O8").Value = 27.20373365408723

# Row 9
$ws.Range("C9").Value = 5.266409101837487
$ws.Range("D9").Value = 8.143925056831105
$ws.Range("E9").Value = 11.8834107526297
$ws.Range("F9").Value = 35.67065886905645
$ws.Range("G9").Value = 3.633253863358509
$ws.Range("J9").Value = 9.192086564531529
$ws.Range("M9").Value = 26.33988729653997
$ws.Range("N9").Value = 17.90827736312335
$ws.Range("O9").Value = 27.70206230748402

# Row 10
$ws.Range("C10").Value = 5.552983637587624
$ws.Range("D10").Value = 8.126088492097622
$ws.Range("E10").Value = 11.79365862909005
$ws.Range("F10").Value = 36.14503534893896
$ws.Range("G10").Value = 3.627431981747816
$ws.Range("J10").Value = 9.12406639034776
$ws.Range("M10").Value = 27.52532825249589
$ws.Range("N10").Value = 18.3747519733078
$ws.Range("O10").Value = 28.10534082559191

# Row 11
$ws.Range("C11").Value = 5.680083447584085
$ws.Range("D11").Value = 8.118488348619831
$ws.Range("E11").Value = 11.75459464164173
$ws.Range("F11").Value = 36.37142326415622
$ws.Range("G11").Value = 3.624906063932233
$ws.Range("J11").Value = 9.094383074525163
$ws.Range("M11").Value = 28.04854415224001
$ws.Range("N11").Value = 18.58480704459731
$ws.Range("O11").Value = 28.2962691859385

# Row 12
$ws.Range("C12").Value = 5.727697019072191
$ws.Range("D12").Value = 8.115683955580511
$ws.Range("E12").Value = 11.74005442393552
$ws.Range("F12").Value = 36.45861144743968
$ws.Range("G12").Value = 3.623967057605458
$ws.Range("J12").Value = 9.083322763094886
$ws.Range("M12").Value = 28.24420637808921
$ws.Range("N12").Value = 18.66397177131545
$ws.Range("O12").Value = 28.36958893509163

# Row 13
$ws.Range("C13").Value = 5.717466240855408
$ws.Range("D13").Value = 8.116284661937106
$ws.Range("E13").Value = 11.7431747120239
$ws.Range("F13").Value = 36.43977001875478
$ws.Range("G13").Value = 3.62416851266289
$ws.Range("J13").Value = 9.085696803982293
$ws.Range("M13").Value = 28.20217918147413
$ws.Range("N13").Value = 18.64694014921102
$ws.Range("O13").Value = 28.35375378897603

# Row 14
$ws.Range("C14").Value = 5.684011243420038
$ws.Range("D14").Value = 8.118256155661193
$ws.Range("E14").Value = 11.75339335677766
$ws.Range("F14").Value = 36.37856740138874
$ws.Range("G14").Value = 3.624828461136186
$ws.Range("J14").Value = 9.093469531882892
$ws.Range("M14").Value = 28.06469173235108
$ws.Range("N14").Value = 18.59132795120682
$ws.Range("O14").Value = 28.30228114332608

# Row 15
$ws.Range("C15").Value = 5.663450543437676
$ws.Range("D15").Value = 8.119473331603901
$ws.Range("E15").Value = 11.75968541198154
$ws.Range("F15").Value = 36.34126723495876
$ws.Range("G15").Value = 3.625234974923587
$ws.Range("J15").Value = 9.098253979413732
$ws.Range("M15").Value = 27.98015076400866
$ws.Range("N15").Value = 18.55721254933057
$ws.Range("O15").Value = 28.27088371984506

# Row 16
$ws.Range("C16").Value = 5.54460754083514
$ws.Range("D16").Value = 8.126595495419462
$ws.Range("E16").Value = 11.79624695199606
$ws.Range("F16").Value = 36.13044818536849
$ws.Range("G16").Value = 3.627599511780836
$ws.Range("J16").Value = 9.12603152301249
$ws.Range("M16").Value = 27.49079784168778
$ws.Range("N16").Value = 18.3609751440959
$ws.Range("O16").Value = 28.09300892633405

# Row 17
$ws.Range("C17").Value = 5.470829981370282
$ws.Range("D17").Value = 8.131096117915977
$ws.Range("E17").Value = 11.81912733171457
$ws.Range("F17").Value = 36.00378655360414
$ws.Range("G17").Value = 3.629081372065449
$ws.Range("J17").Value = 9.143394007956342
$ws.Range("M17").Value = 27.18636605123488
$ws.Range("N17").Value = 18.23998941100385
$ws.Range("O17").Value = 27.98576423220493

# Row 18
$ws.Range("C18").Value = 5.428090272140002
$ws.Range("D18").Value = 8.13373313407673
$ws.Range("E18").Value = 11.83245370037953
$ws.Range("F18").Value = 35.93193508568985
$ws.Range("G18").Value = 3.62994523431616
$ws.Range("J18").Value = 9.15349904969853
$ws.Range("M18").Value = 27.00976304684315
$ws.Range("N18").Value = 18.17020357475227
$ws.Range("O18").Value = 27.92478698553456

# Row 19
$ws.Range("C19").Value = 5.413568509655802
$ws.Range("D19").Value = 8.134634299424102
$ws.Range("E19").Value = 11.83699436085252
$ws.Range("F19").Value = 35.90778119200861
$ws.Range("G19").Value = 3.630239707784114
$ws.Range("J19").Value = 9.156940842831343
$ws.Range("M19").Value = 26.94971551839397
$ws.Range("N19").Value = 18.1465434532122
$ws.Range("O19").Value = 27.90426423104058

# Row 20
$ws.Range("C20").Value = 5.478715632323479
$ws.Range("D20").Value = 8.130612014060203
$ws.Range("E20").Value = 11.81667448670565
$ws.Range("F20").Value = 36.01716675168769
$ws.Range("G20").Value = 3.628922432430944
$ws.Range("J20").Value = 9.141533474329631
$ws.Range("M20").Value = 27.21892999095836
$ws.Range("N20").Value = 18.2528895611365
$ws.Range("O20").Value = 27.99710782964508

# Row 21
$ws.Range("C21").Value = 5.693852148806648
$ws.Range("D21").Value = 8.117675084362947
$ws.Range("E21").Value = 11.75038505153687
$ws.Range("F21").Value = 36.39650498937644
$ws.Range("G21").Value = 3.624634144089717
$ws.Range("J21").Value = 9.091181612815451
$ws.Range("M21").Value = 28.10514331787087
$ws.Range("N21").Value = 18.60767341839836
$ws.Range("O21").Value = 28.31737269687011

# Row 22
$ws.Range("C22").Value = 5.831428445031476
$ws.Range("D22").Value = 8.109649022252045
$ws.Range("E22").Value = 11.70853204454007
$ws.Range("F22").Value = 36.65290286938474
$ws.Range("G22").Value = 3.621933475646634
$ws.Range("J22").Value = 9.059323160145407
$ws.Range("M22").Value = 28.66988462800457
$ws.Range("N22").Value = 18.83730829958182
$ws.Range("O22").Value = 28.53259843002255

# Row 23
$ws.Range("C23").Value = 5.758292576644361
$ws.Range("D23").Value = 8.113893519922838
$ws.Range("E23").Value = 11.73073561469297
$ws.Range("F23").Value = 36.51530441017831
$ws.Range("G23").Value = 3.623365578420934
$ws.Range("J23").Value = 9.076230920216181
$ws.Range("M23").Value = 28.36984271995736
$ws.Range("N23").Value = 18.71497483055848
$ws.Range("O23").Value = 28.41720617782246

# Row 24
$ws.Range("C24").Value = 5.475151535118598
$ws.Range("D24").Value = 8.130830723060848
$ws.Range("E24").Value = 11.81778288167845
$ws.Range("F24").Value = 36.01111454431462
$ws.Range("G24").Value = 3.628994251900087
$ws.Range("J24").Value = 9.142374238023111
$ws.Range("M24").Value = 27.20421274625529
$ws.Range("N24").Value = 18.24705811285196
$ws.Range("O24").Value = 27.99197726963495

# Row 25
$ws.Range("C25").Value = 5.158323954537184
$ws.Range("D25").Value = 8.150950920934912
$ws.Range("E25").Value = 11.91803635589946
$ws.Range("F25").Value = 35.50575169907831
$ws.Range("G25").Value = 3.6355071705681
$ws.Range("J25").Value = 9.218259014512984
$ws.Range("M25").Value = 28.20217918147413
$ws.Range("N25").Value = 18.64694014921102
$ws.Range("O25").Value = 28.35375378897603
